{"js": "// Table 1 caption: \"Table 1: Hypothetical list of AR4D projects\"\n//               -> \"Table 1. Hypothetical list of AR4D projects\"\n// with \"Table 1\", \".\" and \"Hypothetical list of AR4D projects\" all bold,\n// and the separating space left not-bold (its own run), matching the\n// target OOXML (4 runs instead of the original 2).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the caption paragraph robustly (starts with \"Table 1\").\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Table 1\") === 0) {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the 'Table 1' caption paragraph\");\n}\n\nconst whole = target.getRange(\"Whole\");\n\n// Build replacement OOXML: keep the paragraph's own identity/rsid\n// attributes and the first run's rsidRPr, only the run split/text/bold\n// changes per the diff.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  '<w:body>' +\n  '<w:p w14:paraId=\"42B9AC10\" w14:textId=\"77777777\" w:rsidR=\"002E0A54\" w:rsidRDefault=\"00982D68\">' +\n  '<w:r w:rsidRPr=\"00C57713\"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Table 1</w:t></w:r>' +\n  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>.</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Hypothetical list of AR4D projects</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nwhole.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Table 1 caption: \"Table 1: Hypothetical list of AR4D projects\"\n#               -> \"Table 1. Hypothetical list of AR4D projects\"\n# with \"Table 1\", \".\" and \"Hypothetical list of AR4D projects\" all bold,\n# and the separating space left not-bold (its own run), matching the\n# target OOXML (4 runs instead of the original 2).\n\n$d = $word.ActiveDocument\n\n# Locate the caption paragraph robustly via Find, then grab its full\n# paragraph Range so we can rewrite the run layout in one shot.\n$findRng = $d.Content\n$found = $findRng.Find.Execute(\"Table 1:\")\nif (-not $found) {\n    throw \"Could not find the 'Table 1:' caption text\"\n}\n$paraRng = $findRng.Paragraphs(1).Range\n\n$xml = @'\n<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\" w14:paraId=\"42B9AC10\" w14:textId=\"77777777\" w:rsidR=\"002E0A54\" w:rsidRDefault=\"00982D68\"><w:r w:rsidRPr=\"00C57713\"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Table 1</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Hypothetical list of AR4D projects</w:t></w:r></w:p>\n'@\n\n$paraRng.InsertXML($xml)\n"}
